# Applies the cryptos-list data refresh described by the commit message:
#   "Updated cryptos list on Fri Sep 27 13:33:32 UTC 2024 with GitHub Actions"
#
# Every data cell in this sheet (B:E, rows 2-51) is stored as text, including
# the "Price" column (D) which often holds plain-looking numbers such as
# "608.56". Writing such a string straight into `.Value` would make Excel
# auto-detect it as a number and convert the cell, which is not what the
# source data looks like. For any new Price value that would be
# auto-recognised as a number we briefly force the cell format to Text,
# write the value, then clear the format again so the cell keeps its
# original (default) style/look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.653.17'
$ws.Range("E2").Value = '  +2.02%  '

$ws.Range("D3").Value = '2.654.28'
$ws.Range("E3").Value = '  +1.10%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.56'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.81'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.14%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '2.652.93'
$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +9.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.04'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.34%  '

$ws.Range("E12").Value = '  +2.43%  '

$ws.Range("E13").Value = '  +1.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '30.02'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000206'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +20.24%  '

$ws.Range("D16").Value = '3.135.25'
$ws.Range("E16").Value = '  +1.22%  '

$ws.Range("D17").Value = '65.448.61'
$ws.Range("E17").Value = '  +1.78%  '

$ws.Range("D18").Value = '2.652.60'
$ws.Range("E18").Value = '  +2.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.72'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.91'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '359.95'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.71%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.46'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.57%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.33'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.19%  '

$ws.Range("E25").Value = '  +0.52%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.58'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000106'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +16.39%  '

$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.170'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.98%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.23'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +7.71%  '

$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.13'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.28%  '

$ws.Range("E32").Value = '  +0.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '533.79'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.79'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.93%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.53'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.44'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.93%  '

$ws.Range("E37").Value = '  +2.37%  '

$ws.Range("E38").Value = '  +2.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '163.02'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.01'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.35%  '

$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.02'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '166.09'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.17'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.33'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0613'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.15'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.55%  '

$ws.Range("E49").Value = '  +5.84%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.656'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.35%  '

$ws.Range("E51").Value = '  +0.42%  '
